# Capstone TestData workbook: add "Create Claim" / "Cancel Claim" testdata
# sheets (testCreateClaim, testCancelclaim), matching the existing
# 2-row/N-column "data table" layout used by the other sheets in this
# workbook, and touch up a few page/view settings on the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Existing sheets: minor view / print touch-ups
# ---------------------------------------------------------------------

# testValidAdminLogin: selection becomes a full A1:XFD2 band (no longer a
# single active cell)
$ws1 = $wb.Worksheets.Item("testValidAdminLogin")
$ws1.Select()
$ws1.Range("A1:XFD2").Select()
$ws1.PageSetup.Orientation = 1

# testInvalidAdminLogin: just the print setup changes
$ws2 = $wb.Worksheets.Item("testInvalidAdminLogin")
$ws2.PageSetup.Orientation = 1

# testPIMPageLinks: just the print setup changes
$ws3 = $wb.Worksheets.Item("testPIMPageLinks")
$ws3.PageSetup.Orientation = 1

# testCreateUser: just the print setup changes
$ws4 = $wb.Worksheets.Item("testCreateUser")
$ws4.PageSetup.Orientation = 1

# testDashboardSections: print setup changes; tabSelected moves off this
# sheet once the new sheets are added/activated below
$ws5 = $wb.Worksheets.Item("testDashboardSections")
$ws5.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2. New sheets: testCreateClaim, testCancelclaim
#    (copy the existing testCreateUser sheet so the new sheets inherit
#    the same namespaces / row & page formatting as its siblings, then
#    overwrite the data with the claim test data)
# ---------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Copy($null, $lastSheet)
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "testCreateClaim"
$ws6.Cells.Clear()
$ws6.Range("A1").Value = "Username"
$ws6.Range("B1").Value = "Password"
$ws6.Range("C1").Value = "Event_select"
$ws6.Range("D1").Value = "Currency_select"
$ws6.Range("A2").Value = "Admin"
$ws6.Range("B2").Value = "admin123"
$ws6.Range("C2").Value = "Travel Allowance"
$ws6.Range("D2").Value = "Afghanistan Afghani"
$ws6.Range("C1:D2").Select()
$ws6.PageSetup.Orientation = 1

$ws4.Copy($null, $ws6)
$ws7 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7.Name = "testCancelclaim"
$ws7.Cells.Clear()
$ws7.Range("A1").Value = "Username"
$ws7.Range("B1").Value = "Password"
$ws7.Range("C1").Value = "Event_select"
$ws7.Range("D1").Value = "Currency_select"
$ws7.Range("A2").Value = "Admin"
$ws7.Range("B2").Value = "admin123"
$ws7.Range("C2").Value = "Travel Allowance"
$ws7.Range("D2").Value = "Afghanistan Afghani"
$ws7.Range("C1:D2").Select()
$ws7.PageSetup.Orientation = 1

# testCancelclaim ends up the active/selected tab
$ws7.Select()

Write-Host "Applied claim testdata sheets."
